$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Agenda") - bullet "What ABP Team is doing?" -> "AI for ABP Developers"
# Three runs in paragraph 3 of the body textbox (shape 2) get their text
# replaced while keeping their original run-level formatting (bold run in
# the middle stays bold, the others stay regular).
# ---------------------------------------------------------------------------
$slideAgenda = $p.Slides.Item(2)
$bodyAgenda = $slideAgenda.Shapes.Item(2).TextFrame.TextRange

$full = $bodyAgenda.Text
$idx = $full.IndexOf("What ")
$run = $bodyAgenda.Characters($idx + 1, "What ".Length)
$run.Text = "AI for "

$full = $bodyAgenda.Text
$idx = $full.IndexOf("ABP Team ")
$run = $bodyAgenda.Characters($idx + 1, "ABP Team ".Length)
$run.Text = "ABP "

$full = $bodyAgenda.Text
$idx = $full.IndexOf("is doing?")
$run = $bodyAgenda.Characters($idx + 1, "is doing?".Length)
$run.Text = "Developers"

# ---------------------------------------------------------------------------
# Slide 9 ("AI for ABP Developers") - bullet "AI Management Module" gets
# split into "The " (regular) + "AI Management " (bold) + "Module" (regular).
# ---------------------------------------------------------------------------
$slideAI = $p.Slides.Item(9)
$bodyAI = $slideAI.Shapes.Item(2).TextFrame.TextRange

$full = $bodyAI.Text
$idx = $full.IndexOf("AI Management Module")
$run = $bodyAI.Characters($idx + 1, "AI Management Module".Length)
$run.Text = "The AI Management Module"

$full = $bodyAI.Text
$idx = $full.IndexOf("AI Management Module")
$boldRun = $bodyAI.Characters($idx + 1, "AI Management ".Length)
$boldRun.Font.Bold = $true

# ---------------------------------------------------------------------------
# Slide 10 (last slide, "TODO" placeholder) - title run picks up a "dirty"
# proofing flag in the canonical file (no visible/content change).
# ---------------------------------------------------------------------------
$slideLast = $p.Slides.Item(10)
$titleLast = $slideLast.Shapes.Item(1).TextFrame.TextRange
$titleRun = $titleLast.Characters(1, $titleLast.Text.Length)
$titleRun.Text = $titleRun.Text
